# [DSC-844] Improved bitstreams handling during BulkImport
#
# - "bitstream-metadata" sheet: rename the BITSTREAM-ID column header to POSITION
# - The sample bitstream path cell (B2) becomes a real hyperlink ("file://test")
#   whose displayed text is split across two runs: "file://test" (blue) + ".txt"
#   (default colour) instead of the placeholder text that used to be swapped in
#   by the test harness.

$wb = $excel.ActiveWorkbook

# --- bitstream-metadata sheet -------------------------------------------------
$ws4 = $wb.Worksheets.Item("bitstream-metadata")

# Header rename: BITSTREAM-ID -> POSITION
$ws4.Range("D1").Value = "POSITION"

# B2: turn the placeholder text into the real (hyperlinked) file reference,
# with the "file://test" portion coloured blue and ".txt" left default/black.
$b2 = $ws4.Range("B2")

# Hyperlink on B2 pointing at the "file://test" location (this also stamps
# the cell's displayed text, which we immediately re-split into two runs
# below so "file://test" and ".txt" keep their own formatting).
$ws4.Hyperlinks.Add($b2, "file://test", [System.Type]::Missing, [System.Type]::Missing, "file://test")

$b2.Value = "file://test.txt"

$linkPart = $b2.Characters(1, 11)
$linkPart.Font.Name = "Arial"
$linkPart.Font.Size = 10
$linkPart.Font.Color = 16711680
$linkPart.Font.Underline = $false

$extPart = $b2.Characters(12, 4)
$extPart.Font.Name = "Arial"
$extPart.Font.Size = 10
$extPart.Font.ColorIndex = -4105
$extPart.Font.Underline = $false

# Header/footer font used to render "Normale" instead of "Regular"
$ws4.PageSetup.OddHeader = "&C&""Times New Roman,Normale""&12&A"
$ws4.PageSetup.OddFooter = "&C&""Times New Roman,Normale""&12Page &P"

# --- selections on every sheet (mirrors the state Excel leaves behind) -------
$ws1 = $wb.Worksheets.Item("items")
$ws1.Activate()
$ws1.Range("A2").Select()

$ws2 = $wb.Worksheets.Item("dc.contributor.author")
$ws2.Activate()
$ws2.Range("A2").Select()

$ws3 = $wb.Worksheets.Item("dc.contributor.editor")
$ws3.Activate()
$ws3.Range("A2").Select()

$ws4.Activate()
$ws4.Range("B3").Select()
